# Auto-generated Excel COM-interop script applying numeric updates
# to the Leve profit-tracking sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6045.2925
$ws.Range("I62").Value = 3666.3572
$ws.Range("K62").Value = 3666.3572
$ws.Range("M62").Value = -3042.3572
$ws.Range("H65").Value = 6045.2925
$ws.Range("I65").Value = 3666.3572
$ws.Range("K65").Value = 18331.786
$ws.Range("M65").Value = -15211.786
$ws.Range("H76").Value = 5068.4287
$ws.Range("J76").Value = 4996.6665
$ws.Range("L76").Value = 4996.6665
$ws.Range("N76").Value = -5626.6665
$ws.Range("H79").Value = 5068.4287
$ws.Range("J79").Value = 4996.6665
$ws.Range("L79").Value = 4996.6665
$ws.Range("N79").Value = -7180.6665
$ws.Range("H125").Value = 1653.6666
$ws.Range("I125").Value = 732.8
$ws.Range("J125").Value = 2311.4285
$ws.Range("K125").Value = 6595.2
$ws.Range("L125").Value = 20802.8565
$ws.Range("M125").Value = -4135.2
$ws.Range("N125").Value = -25722.8565
$ws.Range("H135").Value = 1463.091
$ws.Range("I135").Value = 1078.8823
$ws.Range("J135").Value = 2769.4
$ws.Range("K135").Value = 9709.940699999999
$ws.Range("L135").Value = 24924.6
$ws.Range("M135").Value = -7174.940699999999
$ws.Range("N135").Value = -29994.6
$ws.Range("H137").Value = 1434.1666
$ws.Range("I137").Value = 1345.85
$ws.Range("J137").Value = 1875.75
$ws.Range("K137").Value = 4037.55
$ws.Range("L137").Value = 5627.25
$ws.Range("M137").Value = -1487.55
$ws.Range("N137").Value = -10727.25
$ws.Range("H138").Value = 3516.468
$ws.Range("I138").Value = 1333.5625
$ws.Range("K138").Value = 4000.6875
$ws.Range("M138").Value = 1139.3125
$ws.Range("H141").Value = 1575.4375
$ws.Range("I141").Value = 1575.4375
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 4726.3125
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 453.6875
$ws.Range("N141").ClearContents()

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3735.7083
$ws.Range("I61").Value = 3267.65
$ws.Range("K61").Value = 3267.65
$ws.Range("M61").Value = -3055.65
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 8729.375
$ws.Range("I74").Value = 1020.86664
$ws.Range("K74").Value = 1020.86664
$ws.Range("M74").Value = -146.86664
$ws.Range("H77").Value = 8729.375
$ws.Range("I77").Value = 1020.86664
$ws.Range("K77").Value = 5104.3332
$ws.Range("M77").Value = -736.3332
$ws.Range("H97").Value = 1236.138
$ws.Range("I97").Value = 1219.75
$ws.Range("J97").Value = 1314.8
$ws.Range("K97").Value = 1219.75
$ws.Range("L97").Value = 1314.8
$ws.Range("M97").Value = -723.75
$ws.Range("N97").Value = -2306.8
$ws.Range("H110").Value = 6765.884
$ws.Range("I110").Value = 8635.166999999999
$ws.Range("J110").Value = 4404.684
$ws.Range("K110").Value = 8635.166999999999
$ws.Range("L110").Value = 4404.684
$ws.Range("M110").Value = -6590.166999999999
$ws.Range("N110").Value = -8494.684000000001
$ws.Range("H122").Value = 2776.1667
$ws.Range("I122").Value = 2106
$ws.Range("K122").Value = 6318
$ws.Range("M122").Value = -3868
$ws.Range("H136").Value = 3735.7083
$ws.Range("I136").Value = 3267.65
$ws.Range("K136").Value = 9802.950000000001
$ws.Range("M136").Value = -7252.950000000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3843.6
$ws.Range("I105").Value = 4433.7144
$ws.Range("J105").Value = 2466.6667
$ws.Range("K105").Value = 4433.7144
$ws.Range("L105").Value = 2466.6667
$ws.Range("M105").Value = -2686.7144
$ws.Range("N105").Value = -5960.6667
$ws.Range("H134").Value = 2111.375
$ws.Range("I134").Value = 1830.74
$ws.Range("K134").Value = 5492.22
$ws.Range("M134").Value = -2957.22

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1355.3334
$ws.Range("I16").Value = 1457
$ws.Range("K16").Value = 1457
$ws.Range("M16").Value = -1170
$ws.Range("H31").Value = 77915.57000000001
$ws.Range("I31").Value = 145399
$ws.Range("J31").Value = 10432.143
$ws.Range("K31").Value = 145399
$ws.Range("L31").Value = 10432.143
$ws.Range("M31").Value = -145104
$ws.Range("N31").Value = -11022.143
$ws.Range("H34").Value = 77915.57000000001
$ws.Range("I34").Value = 145399
$ws.Range("J34").Value = 10432.143
$ws.Range("K34").Value = 145399
$ws.Range("L34").Value = 10432.143
$ws.Range("M34").Value = -145197
$ws.Range("N34").Value = -10836.143
$ws.Range("H107").Value = 276.65216
$ws.Range("I107").Value = 300.58823
$ws.Range("J107").Value = 208.83333
$ws.Range("K107").Value = 300.58823
$ws.Range("L107").Value = 208.83333
$ws.Range("M107").Value = 1619.41177
$ws.Range("N107").Value = -4048.83333
$ws.Range("H113").Value = 1355.3334
$ws.Range("I113").Value = 1457
$ws.Range("K113").Value = 1457
$ws.Range("M113").Value = 713
$ws.Range("H132").Value = 3004.5264
$ws.Range("I132").Value = 2652.1177
$ws.Range("K132").Value = 7956.353099999999
$ws.Range("M132").Value = -5426.353099999999
$ws.Range("H134").Value = 11424.361
$ws.Range("I134").Value = 8202.483
$ws.Range("J134").Value = 31400
$ws.Range("K134").Value = 24607.449
$ws.Range("L134").Value = 94200
$ws.Range("M134").Value = -22072.449
$ws.Range("N134").Value = -99270

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 3939.4546
$ws.Range("I9").Value = 869.8
$ws.Range("J9").Value = 6497.5
$ws.Range("K9").Value = 2609.4
$ws.Range("L9").Value = 19492.5
$ws.Range("M9").Value = -2385.4
$ws.Range("N9").Value = -19940.5
$ws.Range("H14").Value = 134051.73
$ws.Range("I14").Value = 134051.73
$ws.Range("K14").Value = 402155.1900000001
$ws.Range("M14").Value = -401982.1900000001
$ws.Range("H34").Value = 6174525.5
$ws.Range("I34").Value = 43.75
$ws.Range("J34").Value = 11114111
$ws.Range("K34").Value = 131.25
$ws.Range("L34").Value = 33342333
$ws.Range("M34").Value = -47.25
$ws.Range("N34").Value = -33342501
$ws.Range("H39").Value = 1211.1111
$ws.Range("J39").Value = 1350
$ws.Range("L39").Value = 4050
$ws.Range("N39").Value = -4638
$ws.Range("H55").Value = 1810.8
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").ClearContents()
$ws.Range("H121").Value = 17189.834
$ws.Range("I121").Value = 17776.5
$ws.Range("J121").Value = 16896.5
$ws.Range("K121").Value = 53329.5
$ws.Range("L121").Value = 50689.5
$ws.Range("M121").Value = -52019.5
$ws.Range("N121").Value = -53309.5
$ws.Range("H131").Value = 60120.117
$ws.Range("J131").Value = 2432
$ws.Range("L131").Value = 7296
$ws.Range("N131").Value = -17376
$ws.Range("H137").Value = 2758.4546
$ws.Range("J137").Value = 5127.6665
$ws.Range("L137").Value = 15382.9995
$ws.Range("N137").Value = -25582.9995

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1095.5555
$ws.Range("I122").Value = 1101.3334
$ws.Range("K122").Value = 3304.0002
$ws.Range("M122").Value = -854.0001999999999

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3963.76
$ws.Range("I40").Value = 3588.2
$ws.Range("J40").Value = 5466
$ws.Range("K40").Value = 3588.2
$ws.Range("L40").Value = 5466
$ws.Range("M40").Value = -3452.2
$ws.Range("N40").Value = -5738
$ws.Range("H48").Value = 18995
$ws.Range("I48").Value = 18995
$ws.Range("K48").Value = 18995
$ws.Range("M48").Value = -18334
$ws.Range("H132").Value = 3160.8865
$ws.Range("I132").Value = 2533.457
$ws.Range("K132").Value = 7600.370999999999
$ws.Range("M132").Value = -5070.370999999999

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 519.75
$ws.Range("I17").Value = 371.6
$ws.Range("J17").Value = 766.6667
$ws.Range("K17").Value = 371.6
$ws.Range("L17").Value = 766.6667
$ws.Range("M17").Value = -199.6
$ws.Range("N17").Value = -1110.6667
$ws.Range("H122").Value = 3253.1538
$ws.Range("I122").Value = 3175.158
$ws.Range("J122").Value = 3464.8572
$ws.Range("K122").Value = 9525.474
$ws.Range("L122").Value = 10394.5716
$ws.Range("M122").Value = -7075.474
$ws.Range("N122").Value = -15294.5716
$ws.Range("H126").Value = 3098.5833
$ws.Range("I126").Value = 2318.4
$ws.Range("J126").Value = 6999.5
$ws.Range("K126").Value = 6955.200000000001
$ws.Range("L126").Value = 20998.5
$ws.Range("M126").Value = -4485.200000000001
$ws.Range("N126").Value = -25938.5
$ws.Range("H132").Value = 2925.0222
$ws.Range("I132").Value = 3012.919
$ws.Range("J132").Value = 2518.5
$ws.Range("K132").Value = 9038.757
$ws.Range("L132").Value = 7555.5
$ws.Range("M132").Value = -6508.757
$ws.Range("N132").Value = -12615.5
$ws.Range("H136").Value = 2317.6667
$ws.Range("I136").Value = 2261.8
$ws.Range("K136").Value = 6785.400000000001
$ws.Range("M136").Value = -4235.400000000001
